$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 3) mirroring the format of the existing row 2.
$ws.Range("A3").Value = "2024-09-13 19:49:11"
$ws.Range("B3").Value = "MOCK_check_availability"
$ws.Range("C3").Value = "MOCKURL_https://www.opentable.com/r/bar-spero-washington/"
$ws.Range("D3").Value = "MOCK_No availability for the selected date."

# "2024-09-13" looks like a date to Excel's auto-detection, so force it to
# be stored as plain text (matching the source row) using a leading
# apostrophe, then restore the default "Normal" style so no extra
# number-format styling sticks to the cell.
$ws.Range("E3").Value = "'2024-09-13"
$ws.Range("E3").Style = "Normal"

$ws.Range("F3").Value = "19:49:11"
